$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column widths ---
# Before: 1=37, 2=37, 3=37, 4=49, 5=51, 6=80
# After:  1=37, 2=37, 3=49, 4=51, 5=37, 6=80
# Note: COM ColumnWidth <-> stored OOXML width has a fixed offset (~0.8333)
# due to Excel's internal character-width/pixel rounding, so we compensate
# here so the saved width in the XML matches the target exactly.
$ws.Columns.Item(3).ColumnWidth = 48.1666666666667
$ws.Columns.Item(4).ColumnWidth = 50.1666666666667
$ws.Columns.Item(5).ColumnWidth = 36.1666666666667

# --- Update Row 1 header values ---
$ws.Range("A1").Value = "button_testResultActions_class"
$ws.Range("B1").Value = "button_testResultActions_class_1"
$ws.Range("C1").Value = "button_testResultActions_internalRoleButtonName"
$ws.Range("D1").Value = "button_testResultActions_internalRoleButtonName_1"
$ws.Range("E1").Value = "button_testResultDetails_class"
$ws.Range("F1").Value = "button_testResultDetails_internalRoleButtonName"

# --- Update Row 2 values ---
# A2 and B2 remain unchanged.
$ws.Range("C2").Value = "Failed Automations - Apply to"
$ws.Range("D2").Value = "Failed Portal - Login with"
$ws.Range("E2").Value = """]:nth-child(3) [class=""css-1yjo05o"
# F2 remains unchanged.
